$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New single sample/example row replacing the old two example rows
# (app now supports uploading more than one file).
$ws.Range("A2").Value = "Family Office"
$ws.Range("B2").Value = "Pacific Eagle"
$ws.Range("C2").Value = "3rd Party Risk"
$ws.Range("D2").Value = "PJR20242484"
$ws.Range("E2").Value = "SOWVF0202402"
$ws.Range("F2").Value = "Dennis"

# Write "9/26/2024" as a literal text formula result, then flatten it to a
# plain value via copy / paste-values, so the "Created Date" column keeps
# the text "9/26/2024" instead of being auto-converted into a date serial
# number (and without leaving a stray number-format style behind).
$ws.Range("G2").Formula = "=""9/26/2024"""
$ws.Range("G2").Copy()
$ws.Range("G2").PasteSpecial(-4163)

$ws.Range("H2").Value = "Cyber Security"
$ws.Range("I2").Value = $true
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = "In Progress"

# Reset row 2's special styling (centered / vertical-wrap) back to the default,
# then drop the now-unused "Budget Amount" value entirely (no column M value
# for this entry).
$ws.Range("A2:O2").Style = "Normal"
$ws.Range("M2").ClearContents()
$ws.Rows.Item(2).AutoFit()

# Drop the second/old example row entirely.
$ws.Rows.Item(3).Delete()

$ws.Range("O3").Select() | Out-Null
